$wb = $excel.ActiveWorkbook

# Overview sheet: these three source files finished translation and moved
# from "Ready for handoff" to "In Translation" for both zh-cn (col E) and
# de-de (col F) status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"
$wsOverview.Range("E7").Value = "In Translation"
$wsOverview.Range("F7").Value = "In Translation"
$wsOverview.Columns("E:F").AutoFit()

# zh-cn sheet: Status column (C) for the same three rows
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("C6").Value = "In Translation"
$wsZhCn.Range("C7").Value = "In Translation"
$wsZhCn.Columns("C:C").AutoFit()

# de-de sheet: Status column (C) for the same three rows
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("C6").Value = "In Translation"
$wsDeDe.Range("C7").Value = "In Translation"
$wsDeDe.Columns("C:C").AutoFit()
